$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-13 with refreshed TPM-derived values and corrected
# sending/target cluster + ligand/receptor symbol layout.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efemp1"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9023566666666666
$ws.Range("H2").Value = 2.70707
$ws.Range("I2").Value = 0.01513132473647763
$ws.Range("J2").Value = 0.01513132473647763
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 1.237019398558889
$ws.Range("R2").Value = 11.13317458703
$ws.Range("S2").Value = 0.0001669080913850617
$ws.Range("T2").Value = 0.0001669080913850616

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efemp1"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9023566666666666
$ws.Range("H3").Value = 2.70707
$ws.Range("I3").Value = 0.01513132473647763
$ws.Range("J3").Value = 0.01513132473647763
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 83.84131857701222
$ws.Range("R3").Value = 754.5718671931099
$ws.Range("S3").Value = 0.01131251011843356
$ws.Range("T3").Value = 0.01131251011843356

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efemp1"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9023566666666666
$ws.Range("H4").Value = 2.70707
$ws.Range("I4").Value = 0.01513132473647763
$ws.Range("J4").Value = 0.01513132473647763
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 26.81667216062666
$ws.Range("R4").Value = 241.35004944564
$ws.Range("S4").Value = 0.003618309925328173
$ws.Range("T4").Value = 0.003618309925328172

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efemp1"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9023566666666666
$ws.Range("H5").Value = 2.70707
$ws.Range("I5").Value = 0.01513132473647763
$ws.Range("J5").Value = 0.01513132473647763
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.2489972009566666
$ws.Range("R5").Value = 2.24097480861
$ws.Range("S5").Value = 0.00003359660133084114
$ws.Range("T5").Value = 0.00003359660133084113

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efemp1"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 58.36517566666667
$ws.Range("H6").Value = 175.095527
$ws.Range("I6").Value = 0.9787066012115266
$ws.Range("J6").Value = 0.9787066012115266
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 80.01143801227589
$ws.Range("R6").Value = 720.102942110483
$ws.Range("S6").Value = 0.01079575342404575
$ws.Range("T6").Value = 0.01079575342404575

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efemp1"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 58.36517566666667
$ws.Range("H7").Value = 175.095527
$ws.Range("I7").Value = 0.9787066012115266
$ws.Range("J7").Value = 0.9787066012115266
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("Q7").Value = 5422.925842559242
$ws.Range("R7").Value = 48806.33258303317
$ws.Range("S7").Value = 0.7317025126354163
$ws.Range("T7").Value = 0.7317025126354162

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efemp1"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 58.36517566666667
$ws.Range("H8").Value = 175.095527
$ws.Range("I8").Value = 0.9787066012115266
$ws.Range("J8").Value = 0.9787066012115266
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 1734.524539207023
$ws.Range("R8").Value = 15610.7208528632
$ws.Range("S8").Value = 0.234035279185491
$ws.Range("T8").Value = 0.2340352791854909

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efemp1"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 58.36517566666667
$ws.Range("H9").Value = 175.095527
$ws.Range("I9").Value = 0.9787066012115266
$ws.Range("J9").Value = 0.9787066012115266
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 16.10534493863567
$ws.Range("R9").Value = 144.948104447721
$ws.Range("S9").Value = 0.002173055966573651
$ws.Range("T9").Value = 0.00217305596657365

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efemp1"
$ws.Range("C10").Value = "Egfr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3674753333333333
$ws.Range("H10").Value = 1.102426
$ws.Range("I10").Value = 0.006162074051995734
$ws.Range("J10").Value = 0.006162074051995733
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 0.5037632375504444
$ws.Range("R10").Value = 4.533869137953999
$ws.Range("S10").Value = 0.00006797157796188056
$ws.Range("T10").Value = 0.00006797157796188053

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Efemp1"
$ws.Range("C11").Value = "Egfr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3674753333333333
$ws.Range("H11").Value = 1.102426
$ws.Range("I11").Value = 0.006162074051995734
$ws.Range("J11").Value = 0.006162074051995733
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 34.14350182063311
$ws.Range("R11").Value = 307.291516385698
$ws.Range("S11").Value = 0.004606901661140729
$ws.Range("T11").Value = 0.004606901661140728

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efemp1"
$ws.Range("C12").Value = "Egfr"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3674753333333333
$ws.Range("H12").Value = 1.102426
$ws.Range("I12").Value = 0.006162074051995734
$ws.Range("J12").Value = 0.006162074051995733
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 10.92080981406133
$ws.Range("R12").Value = 98.28728832655199
$ws.Range("S12").Value = 0.001473518947696157
$ws.Range("T12").Value = 0.001473518947696157

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efemp1"
$ws.Range("C13").Value = "Egfr"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3674753333333333
$ws.Range("H13").Value = 1.102426
$ws.Range("I13").Value = 0.006162074051995734
$ws.Range("J13").Value = 0.006162074051995733
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 0.1014015109553333
$ws.Range("R13").Value = 0.9126135985979998
$ws.Range("S13").Value = 0.00001368186519696716
$ws.Range("T13").Value = 0.00001368186519696715

# Drop the former rows 14-17 (Resolving-Mac as sending cluster); the sheet
# now only has 12 data rows, shrinking the used range to A1:T13.
$ws.Range("A14:T17").Delete()
